# Updated to new APIX YAML format with multiple APIs support.
# Adds snowData/gateway/consumers/consumingCountryGroups columns and
# repoints the old "description"/"owner_team"/"contact_email" columns to
# the new apiContractUrl / documentationUrl fields, shifting everything
# from column I onward and appending new columns up to Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1) - new column layout A..Q
# ---------------------------------------------------------------------
# Columns K:Q are brand new cells (they don't exist in the source sheet
# yet), so first clone the header style (bold / bordered / centered,
# same as the existing A1:J1 header cells) onto them before writing the
# new header text, keeping the whole header row visually consistent.
$ws.Range("A1").Copy()
$ws.Range("K1:Q1").PasteSpecial(-4122)

$ws.Range("A1").Value = "repository_url"
$ws.Range("B1").Value = "api_technical_name"
$ws.Range("C1").Value = "version"
$ws.Range("D1").Value = "snow_business_application_id"
$ws.Range("E1").Value = "snow_application_service_id"
$ws.Range("F1").Value = "platform"
$ws.Range("G1").Value = "lifecycle_status"
$ws.Range("H1").Value = "classification"
$ws.Range("I1").Value = "api_contract_url"
$ws.Range("J1").Value = "documentation_url"
$ws.Range("K1").Value = "api_hosting_country"
$ws.Range("L1").Value = "gateway_type"
$ws.Range("M1").Value = "gateway_proxy_url"
$ws.Range("N1").Value = "gateway_config_url"
$ws.Range("O1").Value = "consumer_application_service_ids"
$ws.Range("P1").Value = "consuming_country_code"
$ws.Range("Q1").Value = "consuming_group_member_code"

# ---------------------------------------------------------------------
# Row 2 - jenkins-mcp-api replaced with HASE customer-accounts API
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "cto-mmf-hk-hase-pa-customer-accounts"
$ws.Range("C2").Value = "1.0.0"
$ws.Range("D2").Value = "BA0001754"
$ws.Range("E2").Value = "AS0003321"
$ws.Range("F2").Value = "GCP_OTHER"
$ws.Range("G2").Value = "ACTIVE"
$ws.Range("H2").Value = "INTERNAL"
$ws.Range("I2").Value = "https://example.com/api/contract.yaml"
$ws.Range("J2").Value = "https://confluence.example.com/api-docs"
$ws.Range("K2").Value = "GB"
$ws.Range("L2").Value = "KONG"
$ws.Range("M2").Value = "https://proxy.example.com"
$ws.Range("N2").Value = "https://config.example.com"
$ws.Range("O2").Value = "AS0003321"
$ws.Range("P2").Value = "GB"
$ws.Range("Q2").Value = "HASE"

# ---------------------------------------------------------------------
# Row 3 - jenkins-pipeline-api data updated to the new APIX fields
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "2.0.0"
$ws.Range("D3").Value = "BA0001754"
$ws.Range("E3").Value = "AS0003322"
$ws.Range("F3").Value = "GCP_OTHER"
$ws.Range("G3").Value = "ACTIVE"
$ws.Range("H3").Value = "INTERNAL"
$ws.Range("I3").Value = "https://example.com/api/pipeline-contract.yaml"
$ws.Range("J3").Value = "https://confluence.example.com/pipeline-docs"
$ws.Range("K3").Value = "GB"
$ws.Range("L3").Value = "KONG"
$ws.Range("M3").Value = "https://proxy.example.com"
$ws.Range("N3").Value = "https://config.example.com"
$ws.Range("O3").Value = "AS0003322"
$ws.Range("P3").Value = "GB"
$ws.Range("Q3").Value = "HASE"

# ---------------------------------------------------------------------
# Rows 4-12 keep their existing A:J data; only the new K:Q columns are
# introduced (left blank - no snow/gateway/consumer data for these rows).
# ---------------------------------------------------------------------
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("K$r").Value = ""
    $ws.Range("L$r").Value = ""
    $ws.Range("M$r").Value = ""
    $ws.Range("N$r").Value = ""
    $ws.Range("O$r").Value = ""
    $ws.Range("P$r").Value = ""
    $ws.Range("Q$r").Value = ""
}
